$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 1323.25
$ws.Range("I4").Value = 97.666664
$ws.Range("K4").Value = 97.666664
$ws.Range("M4").Value = 16.333336

$ws.Range("H88").Value = 2700
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 2700
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 2700
$ws.Range("M88").Value = ""
$ws.Range("N88").Value = -3512

$ws.Range("H91").Value = 2700
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 2700
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 2700
$ws.Range("M91").Value = ""
$ws.Range("N91").Value = -5508

$ws.Range("H96").Value = 652.7778
$ws.Range("I96").Value = 659.5
$ws.Range("K96").Value = 1978.5
$ws.Range("M96").Value = -605.5

$ws.Range("H106").Value = 4565936.5
$ws.Range("I106").Value = 4803354.5
$ws.Range("K106").Value = 4803354.5
$ws.Range("M106").Value = -4802723.5

$ws.Range("H116").Value = 24927.75
$ws.Range("I116").Value = 50005
$ws.Range("J116").Value = 16568.666
$ws.Range("K116").Value = 50005
$ws.Range("L116").Value = 16568.666
$ws.Range("M116").Value = -46563
$ws.Range("N116").Value = -23452.666

$ws.Range("H134").Value = 44332.668
$ws.Range("J134").Value = 44332.668
$ws.Range("L134").Value = 44332.668
$ws.Range("N134").Value = -54472.668

$ws.Range("H138").Value = 3854.1794
$ws.Range("J138").Value = 3828.6365
$ws.Range("L138").Value = 11485.9095
$ws.Range("N138").Value = -21765.9095

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 70
$ws.Range("I4").Value = 70
$ws.Range("K4").Value = 70
$ws.Range("M4").Value = 46

$ws.Range("H5").Value = 30
$ws.Range("I5").Value = 30
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 30
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 82
$ws.Range("N5").Value = ""

$ws.Range("H32").Value = 5570.7144
$ws.Range("I32").Value = 5472.9473
$ws.Range("K32").Value = 5472.9473
$ws.Range("M32").Value = -5185.9473

$ws.Range("H61").Value = 12504.292
$ws.Range("I61").Value = 10836.5625
$ws.Range("K61").Value = 10836.5625
$ws.Range("M61").Value = -10624.5625

$ws.Range("H132").Value = 2753.5334
$ws.Range("I132").Value = 2753.5334
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 8260.600199999999
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -5730.600199999999
$ws.Range("N132").Value = ""

$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").Value = ""

$ws.Range("H134").Value = 85000
$ws.Range("J134").Value = 85000
$ws.Range("L134").Value = 85000
$ws.Range("N134").Value = -95140

$ws.Range("H135").Value = 79000
$ws.Range("J135").Value = 79000
$ws.Range("L135").Value = 79000
$ws.Range("N135").Value = -89140

$ws.Range("H136").Value = 12504.292
$ws.Range("I136").Value = 10836.5625
$ws.Range("K136").Value = 32509.6875
$ws.Range("M136").Value = -29959.6875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 30
$ws.Range("I4").Value = 30
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 30
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 85
$ws.Range("N4").Value = ""

$ws.Range("H86").Value = 8665
$ws.Range("I86").Value = 1298.8334
$ws.Range("J86").Value = 14189.625
$ws.Range("K86").Value = 1298.8334
$ws.Range("L86").Value = 14189.625
$ws.Range("M86").Value = -175.8334
$ws.Range("N86").Value = -16435.625

$ws.Range("H89").Value = 8665
$ws.Range("I89").Value = 1298.8334
$ws.Range("J89").Value = 14189.625
$ws.Range("K89").Value = 6494.166999999999
$ws.Range("L89").Value = 70948.125
$ws.Range("M89").Value = -878.1669999999995
$ws.Range("N89").Value = -82180.125

$ws.Range("H99").Value = 8605
$ws.Range("I99").Value = 11240.917
$ws.Range("K99").Value = 11240.917
$ws.Range("M99").Value = -9742.916999999999

$ws.Range("H132").Value = 64746.168
$ws.Range("J132").Value = 67695.39999999999
$ws.Range("L132").Value = 67695.39999999999
$ws.Range("N132").Value = -77815.39999999999

$ws.Range("H134").Value = 8365.450999999999
$ws.Range("I134").Value = 9030.317999999999
$ws.Range("K134").Value = 27090.954
$ws.Range("M134").Value = -24555.954

$ws.Range("H139").Value = 250000
$ws.Range("J139").Value = 250000
$ws.Range("L139").Value = 250000
$ws.Range("N139").Value = -260280

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2583.6
$ws.Range("I31").Value = 1042.7333
$ws.Range("K31").Value = 1042.7333
$ws.Range("M31").Value = -747.7333000000001

$ws.Range("H34").Value = 2583.6
$ws.Range("I34").Value = 1042.7333
$ws.Range("K34").Value = 1042.7333
$ws.Range("M34").Value = -840.7333000000001

$ws.Range("H41").Value = 22829.666
$ws.Range("J41").Value = 29245
$ws.Range("L41").Value = 29245
$ws.Range("N41").Value = -30101

$ws.Range("H42").Value = 15000
$ws.Range("J42").Value = 15000
$ws.Range("L42").Value = 15000
$ws.Range("N42").Value = -16186

$ws.Range("H60").Value = 31249.965
$ws.Range("I60").Value = 13999.8
$ws.Range("J60").Value = 35000
$ws.Range("K60").Value = 13999.8
$ws.Range("L60").Value = 35000
$ws.Range("M60").Value = -13488.8
$ws.Range("N60").Value = -36022

$ws.Range("H107").Value = 4644.4165
$ws.Range("J107").Value = 1362
$ws.Range("L107").Value = 1362
$ws.Range("N107").Value = -5202

$ws.Range("H132").Value = 3061.3635
$ws.Range("I132").Value = 3061.3635
$ws.Range("K132").Value = 9184.0905
$ws.Range("M132").Value = -6654.0905

$ws.Range("H134").Value = 6291.9585
$ws.Range("I134").Value = 6147.1875
$ws.Range("J134").Value = 6581.5
$ws.Range("K134").Value = 18441.5625
$ws.Range("L134").Value = 19744.5
$ws.Range("M134").Value = -15906.5625
$ws.Range("N134").Value = -24814.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 7874.5
$ws.Range("I113").Value = 8832.833000000001
$ws.Range("J113").Value = 4999.5
$ws.Range("K113").Value = 8832.833000000001
$ws.Range("L113").Value = 4999.5
$ws.Range("M113").Value = -6662.833000000001
$ws.Range("N113").Value = -9339.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1190.875
$ws.Range("J22").Value = 1203.5333
$ws.Range("L22").Value = 1203.5333
$ws.Range("N22").Value = -1793.5333

$ws.Range("H27").Value = 1190.875
$ws.Range("J27").Value = 1203.5333
$ws.Range("L27").Value = 1203.5333
$ws.Range("N27").Value = -1417.5333

$ws.Range("H40").Value = 2257

$ws.Range("H93").Value = 1596.0588
$ws.Range("I93").Value = 1431.9
$ws.Range("J93").Value = 1830.5714
$ws.Range("K93").Value = 1431.9
$ws.Range("L93").Value = 1830.5714
$ws.Range("M93").Value = -183.9000000000001
$ws.Range("N93").Value = -4326.5714

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 1003635.06
$ws.Range("I14").Value = 1003635.06
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 1003635.06
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -1003467.06
$ws.Range("N14").Value = ""

$ws.Range("H122").Value = 237319.44
$ws.Range("I122").Value = 316710.38
$ws.Range("J122").Value = 6364
$ws.Range("K122").Value = 950131.14
$ws.Range("L122").Value = 19092
$ws.Range("M122").Value = -947681.14
$ws.Range("N122").Value = -23992

$ws.Range("H126").Value = 2940.8
$ws.Range("I126").Value = 3234.6667
$ws.Range("J126").Value = 2500
$ws.Range("K126").Value = 9704.000100000001
$ws.Range("L126").Value = 7500
$ws.Range("M126").Value = -7234.000100000001
$ws.Range("N126").Value = -12440
